$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1349.5
$ws.Range("J40").Value = 1349.5
$ws.Range("L40").Value = 1349.5
$ws.Range("N40").Value = -1699.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3017.5
$ws.Range("J62").Value = 3450
$ws.Range("L62").Value = 3450
$ws.Range("N62").Value = -4698

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3017.5
$ws.Range("J65").Value = 3450
$ws.Range("L65").Value = 17250
$ws.Range("N65").Value = -23490

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1590.7142
$ws.Range("I70").Value = 1316.5
$ws.Range("J70").Value = 1700.4
$ws.Range("K70").Value = 3949.5
$ws.Range("L70").Value = 5101.200000000001
$ws.Range("M70").Value = -3679.5
$ws.Range("N70").Value = -5641.200000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1590.7142
$ws.Range("I73").Value = 1316.5
$ws.Range("J73").Value = 1700.4
$ws.Range("K73").Value = 3949.5
$ws.Range("L73").Value = 5101.200000000001
$ws.Range("M73").Value = -3013.5
$ws.Range("N73").Value = -6973.200000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 220041570
$ws.Range("I86").Value = 275050940
$ws.Range("K86").Value = 275050940
$ws.Range("M86").Value = -275049817

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 220041570
$ws.Range("I89").Value = 275050940
$ws.Range("K89").Value = 1375254700
$ws.Range("M89").Value = -1375249084

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 12501986
$ws.Range("I106").Value = 25001126
$ws.Range("J106").Value = 2846.5
$ws.Range("K106").Value = 25001126
$ws.Range("L106").Value = 2846.5
$ws.Range("M106").Value = -25000495
$ws.Range("N106").Value = -4108.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2430.4092
$ws.Range("I132").Value = 2420.077
$ws.Range("J132").Value = 2511
$ws.Range("K132").Value = 7260.231000000001
$ws.Range("L132").Value = 7533
$ws.Range("M132").Value = -4730.231000000001
$ws.Range("N132").Value = -12593

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2035.2727
$ws.Range("J138").Value = 1978.3556
$ws.Range("L138").Value = 5935.066800000001
$ws.Range("N138").Value = -16215.0668

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 5082.5806
$ws.Range("I141").Value = 2123.75
$ws.Range("K141").Value = 6371.25
$ws.Range("M141").Value = -1191.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 5000
$ws.Range("I22").Value = 5000
$ws.Range("K22").Value = 5000
$ws.Range("M22").Value = -4701

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2996.862
$ws.Range("I61").Value = 2312.2942
$ws.Range("J61").Value = 3966.6667
$ws.Range("K61").Value = 2312.2942
$ws.Range("L61").Value = 3966.6667
$ws.Range("M61").Value = -2100.2942
$ws.Range("N61").Value = -4390.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H68").Value = 53333.332
$ws.Range("J68").Value = 53333.332
$ws.Range("L68").Value = 53333.332
$ws.Range("N68").Value = -54955.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H71").Value = 53333.332
$ws.Range("J71").Value = 53333.332
$ws.Range("L71").Value = 159999.996
$ws.Range("N71").Value = -168111.996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2336.125
$ws.Range("I102").Value = 2336.125
$ws.Range("K102").Value = 2336.125
$ws.Range("M102").Value = -714.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2996.862
$ws.Range("I136").Value = 2312.2942
$ws.Range("J136").Value = 3966.6667
$ws.Range("K136").Value = 6936.882599999999
$ws.Range("L136").Value = 11900.0001
$ws.Range("M136").Value = -4386.882599999999
$ws.Range("N136").Value = -17000.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1302.3077
$ws.Range("I99").Value = 1284.4546
$ws.Range("K99").Value = 1284.4546
$ws.Range("M99").Value = 213.5454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 20835928
$ws.Range("I105").Value = 25002434
$ws.Range("K105").Value = 25002434
$ws.Range("M105").Value = -25000687

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3673.077
$ws.Range("I134").Value = 3659.0908
$ws.Range("J134").Value = 3750
$ws.Range("K134").Value = 10977.2724
$ws.Range("L134").Value = 11250
$ws.Range("M134").Value = -8442.2724
$ws.Range("N134").Value = -16320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5347.357
$ws.Range("I31").Value = 1312.4706
$ws.Range("J31").Value = 11583.091
$ws.Range("K31").Value = 1312.4706
$ws.Range("L31").Value = 11583.091
$ws.Range("M31").Value = -1017.4706
$ws.Range("N31").Value = -12173.091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5347.357
$ws.Range("I34").Value = 1312.4706
$ws.Range("J34").Value = 11583.091
$ws.Range("K34").Value = 1312.4706
$ws.Range("L34").Value = 11583.091
$ws.Range("M34").Value = -1110.4706
$ws.Range("N34").Value = -11987.091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 882.2
$ws.Range("I105").Value = 1000
$ws.Range("K105").Value = 1000
$ws.Range("M105").Value = 747

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1804.08
$ws.Range("I122").Value = 1435.3334
$ws.Range("J122").Value = 1920.5264
$ws.Range("K122").Value = 4306.0002
$ws.Range("L122").Value = 5761.5792
$ws.Range("M122").Value = -1856.0002
$ws.Range("N122").Value = -10661.5792

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2663.6
$ws.Range("I134").Value = 2662.2856
$ws.Range("J134").Value = 2666.6667
$ws.Range("K134").Value = 7986.8568
$ws.Range("L134").Value = 8000.000100000001
$ws.Range("M134").Value = -5451.8568
$ws.Range("N134").Value = -13070.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 12772.7
$ws.Range("J110").Value = 16100
$ws.Range("L110").Value = 48300
$ws.Range("N110").Value = -56480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5396.4443
$ws.Range("J70").Value = 5429.4546
$ws.Range("L70").Value = 5429.4546
$ws.Range("N70").Value = -5969.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5396.4443
$ws.Range("J73").Value = 5429.4546
$ws.Range("L73").Value = 5429.4546
$ws.Range("N73").Value = -7301.4546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2333.5557
$ws.Range("I82").Value = 1833.6666
$ws.Range("J82").Value = 3333.3333
$ws.Range("K82").Value = 1833.6666
$ws.Range("L82").Value = 3333.3333
$ws.Range("M82").Value = -1472.6666
$ws.Range("N82").Value = -4055.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2333.5557
$ws.Range("I85").Value = 1833.6666
$ws.Range("J85").Value = 3333.3333
$ws.Range("K85").Value = 1833.6666
$ws.Range("L85").Value = 3333.3333
$ws.Range("M85").Value = -585.6666
$ws.Range("N85").Value = -5829.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2998.8333
$ws.Range("I100").Value = 2334.3333
$ws.Range("K100").Value = 2334.3333
$ws.Range("M100").Value = -1793.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4124.227
$ws.Range("I136").Value = 4071.182
$ws.Range("J136").Value = 4177.273
$ws.Range("K136").Value = 12213.546
$ws.Range("L136").Value = 12531.819
$ws.Range("M136").Value = -9663.545999999998
$ws.Range("N136").Value = -17631.819
